$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns with refreshed
# data, keeping price cells as plain text (they are not necessarily valid
# numeric literals, e.g. "66.227.54") and without altering cell styling.

# Row 2
$cellD = $ws.Cells.Item(2, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '66.227.54'
$cellD.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +0.65%  '

# Row 3
$cellD = $ws.Cells.Item(3, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '3.320.60'
$cellD.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +0.76%  '

# Row 4
$cellD = $ws.Cells.Item(4, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '0.999'
$cellD.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.20%  '

# Row 5
$cellD = $ws.Cells.Item(5, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '564.04'
$cellD.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.18%  '

# Row 6
$cellD = $ws.Cells.Item(6, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '185.98'
$cellD.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +1.27%  '

# Row 7
$cellD = $ws.Cells.Item(7, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '1.00'
$cellD.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.04%  '

# Row 8
$cellD = $ws.Cells.Item(8, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '3.311.22'
$cellD.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +0.59%  '

# Row 9
$cellD = $ws.Cells.Item(9, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '0.576'
$cellD.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -2.06%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -4.90%  '

# Row 11
$cellD = $ws.Cells.Item(11, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '0.574'
$cellD.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -1.86%  '

# Row 12
$cellD = $ws.Cells.Item(12, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '46.01'
$cellD.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -2.80%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -0.54%  '

# Row 14
$cellD = $ws.Cells.Item(14, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '3.848.52'
$cellD.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.66%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  -2.16%  '

# Row 16
$cellD = $ws.Cells.Item(16, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '587.15'
$cellD.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -8.93%  '

# Row 17
$cellD = $ws.Cells.Item(17, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '66.149.67'
$cellD.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.53%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +0.37%  '

# Row 19
$cellD = $ws.Cells.Item(19, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '3.319.95'
$cellD.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.71%  '

# Row 20
$cellD = $ws.Cells.Item(20, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '17.69'
$cellD.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -2.32%  '

# Row 21
$cellD = $ws.Cells.Item(21, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '10.93'
$cellD.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -4.02%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  -0.91%  '

# Row 23
$cellD = $ws.Cells.Item(23, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '17.95'
$cellD.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -1.97%  '

# Row 24
$cellD = $ws.Cells.Item(24, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '5.01'
$cellD.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +1.87%  '

# Row 25
$cellD = $ws.Cells.Item(25, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '98.36'
$cellD.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -8.81%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.05%  '

# Row 27
$cellD = $ws.Cells.Item(27, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '2.70'
$cellD.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +0.50%  '

# Row 28
$cellD = $ws.Cells.Item(28, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '9.39'
$cellD.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -2.03%  '

# Row 29
$cellD = $ws.Cells.Item(29, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '8.45'
$cellD.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -2.75%  '

# Row 30
$cellD = $ws.Cells.Item(30, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '30.58'
$cellD.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +0.86%  '

# Row 31
$cellD = $ws.Cells.Item(31, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '6.63'
$cellD.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +5.26%  '

# Row 32
$cellD = $ws.Cells.Item(32, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '3.68'
$cellD.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -5.91%  '

# Row 33
$cellD = $ws.Cells.Item(33, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '560.44'
$cellD.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +7.62%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -2.11%  '

# Row 35
$cellD = $ws.Cells.Item(35, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '3.803.00'
$cellD.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +0.02%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -1.64%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -0.03%  '

# Row 38
$cellD = $ws.Cells.Item(38, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '56.01'
$cellD.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -2.52%  '

# Row 39
$cellD = $ws.Cells.Item(39, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '33.22'
$cellD.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +0.68%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -6.87%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -6.95%  '

# Row 43
$cellD = $ws.Cells.Item(43, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '3.38'
$cellD.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +3.83%  '

# Row 44
$cellD = $ws.Cells.Item(44, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '2.59'
$cellD.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -4.82%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -1.71%  '

# Row 46
$cellD = $ws.Cells.Item(46, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '0.0411'
$cellD.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -0.87%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -9.27%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  -2.44%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +0.06%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -3.03%  '

# Row 51
$cellD = $ws.Cells.Item(51, 4)
$cellD.NumberFormat = "@"
$cellD.Value = '128.83'
$cellD.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +5.63%  '
